$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.251.99"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.835.40"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9985"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.74"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6254"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07392"
$ws.Range("E8").Value = "  -1.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2928"
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.24"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07681"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.821.59"
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.975"
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6703"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.75"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008986"
$ws.Range("E16").Value = "  -3.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.894"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.217.00"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.072.52"
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.77"
$ws.Range("E20").Value = "  +2.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.51"
$ws.Range("E21").Value = "  -1.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.392"
$ws.Range("E23").Value = "  +2.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.43"
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1412"
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.555"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.70"
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.489"
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05820"
$ws.Range("E30").Value = "  +4.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.109"
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.092"
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.208"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7330"
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.146"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.607"
$ws.Range("E37").Value = "  -2.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.848"
$ws.Range("E38").Value = "  +2.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.227.51"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01762"
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.292"
$ws.Range("E41").Value = "  -4.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9161"
$ws.Range("E42").Value = "  +1.71%  "
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.94"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.973.18"
$ws.Range("E45").Value = "  -2.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.29"
$ws.Range("E46").Value = "  -1.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5045"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.182"
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000118"
$ws.Range("E49").Value = "  -3.71%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4039"
$ws.Range("E50").Value = "  -1.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1136"
$ws.Range("E51").Value = "  +3.04%  "
